# Updated symbol list on Sat Dec 17 10:32:56 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-less approach: all "Price" column (D) values are stored as text in the
# original workbook, so force text number-format before writing so Excel does
# not silently convert the numeric-looking strings into real numbers.

function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
}

# --- Price (column D) updates ---
Set-TextValue "D2"  "237.16"
Set-TextValue "D3"  "21.91"
Set-TextValue "D4"  "5.432"
Set-TextValue "D5"  "0.05632"
Set-TextValue "D6"  "6.472"
Set-TextValue "D7"  "3.350"
Set-TextValue "D8"  "1.080"
Set-TextValue "D9"  "0.7886"
Set-TextValue "D10" "0.1398"
Set-TextValue "D11" "0.07335"
Set-TextValue "D12" "0.03208"
Set-TextValue "D13" "0.02973"
Set-TextValue "D14" "0.09247"
Set-TextValue "D15" "0.001670"
Set-TextValue "D16" "3.252"
Set-TextValue "D17" "0.04764"

# Row 18 - Volume(1h) text tweak only
Set-TextValue "E18" "17OneONEWorstin24h"

Set-TextValue "D19" "0.006261"
Set-TextValue "D21" "0.001052"
Set-TextValue "D23" "3.874"
Set-TextValue "D25" "0.3321"
Set-TextValue "D27" "0.0004012"
Set-TextValue "D40" "0.04119"
Set-TextValue "D41" "0.006976"

# Rows 42 and 43 effectively swap coin identity (CEJI <-> BKEXToken) along
# with refreshed price data.
Set-TextValue "B42" "BKEXToken"
Set-TextValue "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1040"
Set-TextValue "E42" "41BKEXTokenBKK"

Set-TextValue "B43" "CEJI"
Set-TextValue "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003041"
Set-TextValue "E43" "42CEJICEJI"

Set-TextValue "D44" "0.009920"
Set-TextValue "D45" "0.00005440"

Set-TextValue "D48" "0.03807"
Set-TextValue "E48" "47BOLOBOLO"

Write-Host "Applied cryptos.xlsx price/coin updates"
